$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C16").Value = '9178229'
$ws.Range("D16").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E16").Value = '1906'
$ws.Range("F16").Value = 128000
$ws.Range("G16").Value = 3200000

$ws.Range("C17").Value = '92504619'
$ws.Range("D17").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E17").Value = '1906'
$ws.Range("F17").Value = 80000
$ws.Range("G17").Value = 2000000

$ws.Range("C18").Value = '9178229'
$ws.Range("D18").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E18").Value = '1907'
$ws.Range("F18").Value = 128000
$ws.Range("G18").Value = 3200000

$ws.Range("C19").Value = '92504619'
$ws.Range("D19").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E19").Value = '1907'
$ws.Range("F19").Value = 80000
$ws.Range("G19").Value = 2000000

$ws.Range("C20").Value = '9178229'
$ws.Range("D20").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E20").Value = '1908'
$ws.Range("F20").Value = 128000
$ws.Range("G20").Value = 3200000

$ws.Range("C21").Value = '92504619'
$ws.Range("D21").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E21").Value = '1908'
$ws.Range("F21").Value = 80000
$ws.Range("G21").Value = 2000000

$ws.Range("C22").Value = '9178229'
$ws.Range("D22").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E22").Value = '1909'
$ws.Range("F22").Value = 128000
$ws.Range("G22").Value = 3200000

$ws.Range("C23").Value = '72133110'
$ws.Range("D23").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E23").Value = '1909'
$ws.Range("F23").Value = 200000
$ws.Range("G23").Value = 5000000

$ws.Range("C24").Value = '92504619'
$ws.Range("D24").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E24").Value = '1909'
$ws.Range("F24").Value = 80000
$ws.Range("G24").Value = 2000000

$ws.Range("C25").Value = '9178229'
$ws.Range("D25").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E25").Value = '1910'
$ws.Range("F25").Value = 128000
$ws.Range("G25").Value = 3200000

$ws.Range("C26").Value = '72133110'
$ws.Range("D26").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E26").Value = '1910'
$ws.Range("F26").Value = 200000
$ws.Range("G26").Value = 5000000

$ws.Range("C27").Value = '92504619'
$ws.Range("D27").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E27").Value = '1910'
$ws.Range("F27").Value = 80000
$ws.Range("G27").Value = 2000000

$ws.Range("C28").Value = '9178229'
$ws.Range("D28").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E28").Value = '1911'
$ws.Range("F28").Value = 128000
$ws.Range("G28").Value = 3200000

$ws.Range("C29").Value = '72133110'
$ws.Range("D29").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E29").Value = '1911'
$ws.Range("F29").Value = 200000
$ws.Range("G29").Value = 5000000

$ws.Range("C30").Value = '92504619'
$ws.Range("D30").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E30").Value = '1911'
$ws.Range("F30").Value = 80000
$ws.Range("G30").Value = 2000000

$ws.Range("C31").Value = '9178229'
$ws.Range("D31").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E31").Value = '1912'
$ws.Range("F31").Value = 128000
$ws.Range("G31").Value = 3200000

$ws.Range("C32").Value = '72133110'
$ws.Range("D32").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E32").Value = '1912'
$ws.Range("F32").Value = 200000
$ws.Range("G32").Value = 5000000

$ws.Range("C33").Value = '92504619'
$ws.Range("D33").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E33").Value = '1912'
$ws.Range("F33").Value = 80000
$ws.Range("G33").Value = 2000000

$ws.Range("C34").Value = '9178229'
$ws.Range("D34").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E34").Value = '2001'
$ws.Range("F34").Value = 128000
$ws.Range("G34").Value = 3200000

$ws.Range("C35").Value = '72133110'
$ws.Range("D35").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E35").Value = '2001'
$ws.Range("F35").Value = 200000
$ws.Range("G35").Value = 5000000

$ws.Range("C36").Value = '92504619'
$ws.Range("D36").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E36").Value = '2001'
$ws.Range("F36").Value = 80000
$ws.Range("G36").Value = 2000000

$ws.Range("C37").Value = '9178229'
$ws.Range("D37").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E37").Value = '2002'
$ws.Range("F37").Value = 128000
$ws.Range("G37").Value = 3200000

$ws.Range("C38").Value = '72133110'
$ws.Range("D38").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E38").Value = '2002'
$ws.Range("F38").Value = 200000
$ws.Range("G38").Value = 5000000

$ws.Range("C39").Value = '92504619'
$ws.Range("D39").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E39").Value = '2002'
$ws.Range("F39").Value = 80000
$ws.Range("G39").Value = 2000000

$ws.Range("C40").Value = '9178229'
$ws.Range("D40").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E40").Value = '2003'
$ws.Range("F40").Value = 128000
$ws.Range("G40").Value = 3200000

$ws.Range("C41").Value = '72133110'
$ws.Range("D41").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E41").Value = '2003'
$ws.Range("F41").Value = 200000
$ws.Range("G41").Value = 5000000

$ws.Range("C42").Value = '92504619'
$ws.Range("D42").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E42").Value = '2003'
$ws.Range("F42").Value = 80000
$ws.Range("G42").Value = 2000000

$ws.Range("C43").Value = '9178229'
$ws.Range("D43").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E43").Value = '2004'
$ws.Range("F43").Value = 128000
$ws.Range("G43").Value = 3200000

$ws.Range("C44").Value = '72133110'
$ws.Range("D44").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E44").Value = '2004'
$ws.Range("F44").Value = 200000
$ws.Range("G44").Value = 5000000

$ws.Range("C45").Value = '92504619'
$ws.Range("D45").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E45").Value = '2004'
$ws.Range("F45").Value = 80000
$ws.Range("G45").Value = 2000000

$ws.Range("C46").Value = '9178229'
$ws.Range("D46").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E46").Value = '2005'
$ws.Range("F46").Value = 128000
$ws.Range("G46").Value = 3200000

$ws.Range("C47").Value = '72133110'
$ws.Range("D47").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E47").Value = '2005'
$ws.Range("F47").Value = 200000
$ws.Range("G47").Value = 5000000

$ws.Range("C48").Value = '92504619'
$ws.Range("D48").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E48").Value = '2005'
$ws.Range("F48").Value = 80000
$ws.Range("G48").Value = 2000000

$ws.Range("C49").Value = '9178229'
$ws.Range("D49").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E49").Value = '2006'
$ws.Range("F49").Value = 128000
$ws.Range("G49").Value = 3200000

$ws.Range("C50").Value = '72133110'
$ws.Range("D50").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E50").Value = '2006'
$ws.Range("F50").Value = 200000
$ws.Range("G50").Value = 5000000

$ws.Range("C51").Value = '92504619'
$ws.Range("D51").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E51").Value = '2006'
$ws.Range("F51").Value = 80000
$ws.Range("G51").Value = 2000000

$ws.Range("C52").Value = '9178229'
$ws.Range("D52").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E52").Value = '2007'
$ws.Range("F52").Value = 128000
$ws.Range("G52").Value = 3200000

$ws.Range("C53").Value = '72133110'
$ws.Range("D53").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E53").Value = '2007'
$ws.Range("F53").Value = 200000
$ws.Range("G53").Value = 5000000

$ws.Range("C54").Value = '92504619'
$ws.Range("D54").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E54").Value = '2007'
$ws.Range("F54").Value = 80000
$ws.Range("G54").Value = 2000000

$ws.Range("C55").Value = '9178229'
$ws.Range("D55").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E55").Value = '2008'
$ws.Range("F55").Value = 128000
$ws.Range("G55").Value = 3200000

$ws.Range("C56").Value = '72133110'
$ws.Range("D56").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E56").Value = '2008'
$ws.Range("F56").Value = 200000
$ws.Range("G56").Value = 5000000

$ws.Range("C57").Value = '92504619'
$ws.Range("D57").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E57").Value = '2008'
$ws.Range("F57").Value = 80000
$ws.Range("G57").Value = 2000000

$ws.Range("C58").Value = '9178229'
$ws.Range("D58").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E58").Value = '2009'
$ws.Range("F58").Value = 128000
$ws.Range("G58").Value = 3200000

$ws.Range("C59").Value = '72133110'
$ws.Range("D59").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E59").Value = '2009'
$ws.Range("F59").Value = 200000
$ws.Range("G59").Value = 5000000

$ws.Range("C60").Value = '92504619'
$ws.Range("D60").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E60").Value = '2009'
$ws.Range("F60").Value = 80000
$ws.Range("G60").Value = 2000000

$ws.Range("C61").Value = '9178229'
$ws.Range("D61").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E61").Value = '2010'
$ws.Range("F61").Value = 128000
$ws.Range("G61").Value = 3200000

$ws.Range("C62").Value = '72133110'
$ws.Range("D62").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E62").Value = '2010'
$ws.Range("F62").Value = 200000
$ws.Range("G62").Value = 5000000

$ws.Range("C63").Value = '92504619'
$ws.Range("D63").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E63").Value = '2010'
$ws.Range("F63").Value = 80000
$ws.Range("G63").Value = 2000000

$ws.Range("C64").Value = '9178229'
$ws.Range("D64").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E64").Value = '2011'
$ws.Range("F64").Value = 128000
$ws.Range("G64").Value = 3200000

$ws.Range("C65").Value = '72133110'
$ws.Range("D65").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E65").Value = '2011'
$ws.Range("F65").Value = 200000
$ws.Range("G65").Value = 5000000

$ws.Range("C66").Value = '92504619'
$ws.Range("D66").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E66").Value = '2011'
$ws.Range("F66").Value = 80000
$ws.Range("G66").Value = 2000000

$ws.Range("C67").Value = '9178229'
$ws.Range("D67").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E67").Value = '2012'
$ws.Range("F67").Value = 128000
$ws.Range("G67").Value = 3200000

$ws.Range("C68").Value = '72133110'
$ws.Range("D68").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E68").Value = '2012'
$ws.Range("F68").Value = 200000
$ws.Range("G68").Value = 5000000

$ws.Range("C69").Value = '92504619'
$ws.Range("D69").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E69").Value = '2012'
$ws.Range("F69").Value = 80000
$ws.Range("G69").Value = 2000000

$ws.Range("C70").Value = '9178229'
$ws.Range("D70").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E70").Value = '2101'
$ws.Range("F70").Value = 128000
$ws.Range("G70").Value = 3200000

$ws.Range("C71").Value = '72133110'
$ws.Range("D71").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E71").Value = '2101'
$ws.Range("F71").Value = 200000
$ws.Range("G71").Value = 5000000

$ws.Range("C72").Value = '92504619'
$ws.Range("D72").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E72").Value = '2101'
$ws.Range("F72").Value = 80000
$ws.Range("G72").Value = 2000000

$ws.Range("C73").Value = '9178229'
$ws.Range("D73").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E73").Value = '2102'
$ws.Range("F73").Value = 128000
$ws.Range("G73").Value = 3200000

$ws.Range("C74").Value = '72133110'
$ws.Range("D74").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E74").Value = '2102'
$ws.Range("F74").Value = 200000
$ws.Range("G74").Value = 5000000

$ws.Range("C75").Value = '92504619'
$ws.Range("D75").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E75").Value = '2102'
$ws.Range("F75").Value = 80000
$ws.Range("G75").Value = 2000000

$ws.Range("C76").Value = '9178229'
$ws.Range("D76").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E76").Value = '2103'
$ws.Range("F76").Value = 128000
$ws.Range("G76").Value = 3200000

$ws.Range("C77").Value = '72133110'
$ws.Range("D77").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E77").Value = '2103'
$ws.Range("F77").Value = 200000
$ws.Range("G77").Value = 5000000

$ws.Range("C78").Value = '92504619'
$ws.Range("D78").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E78").Value = '2103'
$ws.Range("F78").Value = 80000
$ws.Range("G78").Value = 2000000

$ws.Range("C79").Value = '9178229'
$ws.Range("D79").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E79").Value = '2104'
$ws.Range("F79").Value = 128000
$ws.Range("G79").Value = 3200000

$ws.Range("C80").Value = '72133110'
$ws.Range("D80").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E80").Value = '2104'
$ws.Range("F80").Value = 200000
$ws.Range("G80").Value = 5000000

$ws.Range("C81").Value = '92504619'
$ws.Range("D81").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E81").Value = '2104'
$ws.Range("F81").Value = 80000
$ws.Range("G81").Value = 2000000

$ws.Range("C82").Value = '9178229'
$ws.Range("D82").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E82").Value = '2105'
$ws.Range("F82").Value = 128000
$ws.Range("G82").Value = 3200000

$ws.Range("C83").Value = '72133110'
$ws.Range("D83").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E83").Value = '2105'
$ws.Range("F83").Value = 200000
$ws.Range("G83").Value = 5000000

$ws.Range("C84").Value = '92504619'
$ws.Range("D84").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E84").Value = '2105'
$ws.Range("F84").Value = 80000
$ws.Range("G84").Value = 2000000

$ws.Range("C85").Value = '9178229'
$ws.Range("D85").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E85").Value = '2106'
$ws.Range("F85").Value = 128000
$ws.Range("G85").Value = 3200000

$ws.Range("C86").Value = '72133110'
$ws.Range("D86").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E86").Value = '2106'
$ws.Range("F86").Value = 200000
$ws.Range("G86").Value = 5000000

$ws.Range("C87").Value = '92504619'
$ws.Range("D87").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E87").Value = '2106'
$ws.Range("F87").Value = 80000
$ws.Range("G87").Value = 2000000

$ws.Range("C88").Value = '9178229'
$ws.Range("D88").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E88").Value = '2107'
$ws.Range("F88").Value = 128000
$ws.Range("G88").Value = 3200000

$ws.Range("C89").Value = '72133110'
$ws.Range("D89").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E89").Value = '2107'
$ws.Range("F89").Value = 200000
$ws.Range("G89").Value = 5000000

$ws.Range("C90").Value = '92504619'
$ws.Range("D90").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E90").Value = '2107'
$ws.Range("F90").Value = 80000
$ws.Range("G90").Value = 2000000

$ws.Range("C91").Value = '9178229'
$ws.Range("D91").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E91").Value = '2108'
$ws.Range("F91").Value = 128000
$ws.Range("G91").Value = 3200000

$ws.Range("C92").Value = '72133110'
$ws.Range("D92").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E92").Value = '2108'
$ws.Range("F92").Value = 200000
$ws.Range("G92").Value = 5000000

$ws.Range("C93").Value = '92504619'
$ws.Range("D93").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E93").Value = '2108'
$ws.Range("F93").Value = 80000
$ws.Range("G93").Value = 2000000

$ws.Range("C94").Value = '9178229'
$ws.Range("D94").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E94").Value = '2109'
$ws.Range("F94").Value = 128000
$ws.Range("G94").Value = 3200000

$ws.Range("C95").Value = '72133110'
$ws.Range("D95").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E95").Value = '2109'
$ws.Range("F95").Value = 200000
$ws.Range("G95").Value = 5000000

$ws.Range("C96").Value = '92504619'
$ws.Range("D96").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E96").Value = '2109'
$ws.Range("F96").Value = 80000
$ws.Range("G96").Value = 2000000

$ws.Range("C97").Value = '9178229'
$ws.Range("D97").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E97").Value = '2110'
$ws.Range("F97").Value = 128000
$ws.Range("G97").Value = 3200000

$ws.Range("C98").Value = '72133110'
$ws.Range("D98").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E98").Value = '2110'
$ws.Range("F98").Value = 200000
$ws.Range("G98").Value = 5000000

$ws.Range("C99").Value = '92504619'
$ws.Range("D99").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E99").Value = '2110'
$ws.Range("F99").Value = 80000
$ws.Range("G99").Value = 2000000

$ws.Range("C100").Value = '9178229'
$ws.Range("D100").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E100").Value = '2111'
$ws.Range("F100").Value = 128000
$ws.Range("G100").Value = 3200000

$ws.Range("C101").Value = '72133110'
$ws.Range("D101").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E101").Value = '2111'
$ws.Range("F101").Value = 200000
$ws.Range("G101").Value = 5000000

$ws.Range("C102").Value = '92504619'
$ws.Range("D102").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E102").Value = '2111'
$ws.Range("F102").Value = 80000
$ws.Range("G102").Value = 2000000

$ws.Range("C103").Value = '9178229'
$ws.Range("D103").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E103").Value = '2112'
$ws.Range("F103").Value = 128000
$ws.Range("G103").Value = 3200000

$ws.Range("C104").Value = '72133110'
$ws.Range("D104").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E104").Value = '2112'
$ws.Range("F104").Value = 200000
$ws.Range("G104").Value = 5000000

$ws.Range("C105").Value = '92504619'
$ws.Range("D105").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E105").Value = '2112'
$ws.Range("F105").Value = 80000
$ws.Range("G105").Value = 2000000

$ws.Range("C106").Value = '9178229'
$ws.Range("D106").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E106").Value = '2201'
$ws.Range("F106").Value = 128000
$ws.Range("G106").Value = 3200000

$ws.Range("C107").Value = '72133110'
$ws.Range("D107").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E107").Value = '2201'
$ws.Range("F107").Value = 200000
$ws.Range("G107").Value = 5000000

$ws.Range("C108").Value = '92504619'
$ws.Range("D108").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E108").Value = '2201'
$ws.Range("F108").Value = 80000
$ws.Range("G108").Value = 2000000

$ws.Range("C109").Value = '9178229'
$ws.Range("D109").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E109").Value = '2202'
$ws.Range("F109").Value = 128000
$ws.Range("G109").Value = 3200000

$ws.Range("C110").Value = '72133110'
$ws.Range("D110").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E110").Value = '2202'
$ws.Range("F110").Value = 200000
$ws.Range("G110").Value = 5000000

$ws.Range("C111").Value = '92504619'
$ws.Range("D111").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E111").Value = '2202'
$ws.Range("F111").Value = 80000
$ws.Range("G111").Value = 2000000

$ws.Range("C112").Value = '9178229'
$ws.Range("D112").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E112").Value = '2203'
$ws.Range("F112").Value = 128000
$ws.Range("G112").Value = 3200000

$ws.Range("C113").Value = '72133110'
$ws.Range("D113").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E113").Value = '2203'
$ws.Range("F113").Value = 200000
$ws.Range("G113").Value = 5000000

$ws.Range("C114").Value = '92504619'
$ws.Range("D114").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E114").Value = '2203'
$ws.Range("F114").Value = 80000
$ws.Range("G114").Value = 2000000

$ws.Range("C115").Value = '9178229'
$ws.Range("D115").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E115").Value = '2204'
$ws.Range("F115").Value = 128000
$ws.Range("G115").Value = 3200000

$ws.Range("C116").Value = '72133110'
$ws.Range("D116").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E116").Value = '2204'
$ws.Range("F116").Value = 200000
$ws.Range("G116").Value = 5000000

$ws.Range("C117").Value = '92504619'
$ws.Range("D117").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E117").Value = '2204'
$ws.Range("F117").Value = 80000
$ws.Range("G117").Value = 2000000

$ws.Range("C118").Value = '9178229'
$ws.Range("D118").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E118").Value = '2205'
$ws.Range("F118").Value = 128000
$ws.Range("G118").Value = 3200000

$ws.Range("C119").Value = '72133110'
$ws.Range("D119").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E119").Value = '2205'
$ws.Range("F119").Value = 200000
$ws.Range("G119").Value = 5000000

$ws.Range("C120").Value = '92504619'
$ws.Range("D120").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E120").Value = '2205'
$ws.Range("F120").Value = 80000
$ws.Range("G120").Value = 2000000

$ws.Range("C121").Value = '9178229'
$ws.Range("D121").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E121").Value = '2206'
$ws.Range("F121").Value = 128000
$ws.Range("G121").Value = 3200000

$ws.Range("C122").Value = '72133110'
$ws.Range("D122").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E122").Value = '2206'
$ws.Range("F122").Value = 200000
$ws.Range("G122").Value = 5000000

$ws.Range("C123").Value = '92504619'
$ws.Range("D123").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E123").Value = '2206'
$ws.Range("F123").Value = 80000
$ws.Range("G123").Value = 2000000

$ws.Range("C124").Value = '9178229'
$ws.Range("D124").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E124").Value = '2207'
$ws.Range("F124").Value = 128000
$ws.Range("G124").Value = 3200000

$ws.Range("C125").Value = '72133110'
$ws.Range("D125").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E125").Value = '2207'
$ws.Range("F125").Value = 200000
$ws.Range("G125").Value = 5000000

$ws.Range("C126").Value = '92504619'
$ws.Range("D126").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E126").Value = '2207'
$ws.Range("F126").Value = 80000
$ws.Range("G126").Value = 2000000

$ws.Range("C127").Value = '9178229'
$ws.Range("D127").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E127").Value = '2208'
$ws.Range("F127").Value = 128000
$ws.Range("G127").Value = 3200000

$ws.Range("C128").Value = '72133110'
$ws.Range("D128").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E128").Value = '2208'
$ws.Range("F128").Value = 200000
$ws.Range("G128").Value = 5000000

$ws.Range("C129").Value = '92504619'
$ws.Range("D129").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E129").Value = '2208'
$ws.Range("F129").Value = 80000
$ws.Range("G129").Value = 2000000

$ws.Range("C130").Value = '9178229'
$ws.Range("D130").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E130").Value = '2209'
$ws.Range("F130").Value = 128000
$ws.Range("G130").Value = 3200000

$ws.Range("C131").Value = '72133110'
$ws.Range("D131").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E131").Value = '2209'
$ws.Range("F131").Value = 200000
$ws.Range("G131").Value = 5000000

$ws.Range("C132").Value = '92504619'
$ws.Range("D132").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E132").Value = '2209'
$ws.Range("F132").Value = 80000
$ws.Range("G132").Value = 2000000

$ws.Range("C133").Value = '9178229'
$ws.Range("D133").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E133").Value = '2210'
$ws.Range("F133").Value = 128000
$ws.Range("G133").Value = 3200000

$ws.Range("C134").Value = '72133110'
$ws.Range("D134").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E134").Value = '2210'
$ws.Range("F134").Value = 200000
$ws.Range("G134").Value = 5000000

$ws.Range("C135").Value = '92504619'
$ws.Range("D135").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E135").Value = '2210'
$ws.Range("F135").Value = 80000
$ws.Range("G135").Value = 2000000

$ws.Range("C136").Value = '9178229'
$ws.Range("D136").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E136").Value = '2211'
$ws.Range("F136").Value = 128000
$ws.Range("G136").Value = 3200000

$ws.Range("C137").Value = '72133110'
$ws.Range("D137").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E137").Value = '2211'
$ws.Range("F137").Value = 200000
$ws.Range("G137").Value = 5000000

$ws.Range("C138").Value = '92504619'
$ws.Range("D138").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E138").Value = '2211'
$ws.Range("F138").Value = 80000
$ws.Range("G138").Value = 2000000

$ws.Range("C139").Value = '9178229'
$ws.Range("D139").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E139").Value = '2212'
$ws.Range("F139").Value = 128000
$ws.Range("G139").Value = 3200000

$ws.Range("C140").Value = '72133110'
$ws.Range("D140").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E140").Value = '2212'
$ws.Range("F140").Value = 200000
$ws.Range("G140").Value = 5000000

$ws.Range("C141").Value = '92504619'
$ws.Range("D141").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E141").Value = '2212'
$ws.Range("F141").Value = 80000
$ws.Range("G141").Value = 2000000

$ws.Range("C142").Value = '9178229'
$ws.Range("D142").Value = 'CESAR AUGUSTO CARO BARRAZA'
$ws.Range("E142").Value = '2301'
$ws.Range("F142").Value = 106667
$ws.Range("G142").Value = 3200000

$ws.Range("C143").Value = '72133110'
$ws.Range("D143").Value = 'ALVARO JOSE ESMERAL PERTUZ'
$ws.Range("E143").Value = '2301'
$ws.Range("F143").Value = 166667
$ws.Range("G143").Value = 5000000

$ws.Range("C144").Value = '92504619'
$ws.Range("D144").Value = 'HERNANDO DE LA CRUZ VIAÑA'
$ws.Range("E144").Value = '2301'
$ws.Range("F144").Value = 66667
$ws.Range("G144").Value = 2000000
